$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(14, 'Western Michigan', 'Bowling Green', 0, 24, 0, 1),
    @(14, 'Buffalo', 'Eastern Michigan', 0, 13, 0, 1),
    @(14, 'Mississippi State', 'Mississippi', 0, 10, 0, 1),
    @(14, 'Boise State', 'Air Force', 0, 8, 1, 0),
    @(14, 'Nebraska', 'Iowa', 0, 3, 0, 1),
    @(14, 'Temple', 'Memphis', 0, 24, 0, 1),
    @(14, 'Boston College', 'Miami (FL)', 0, 25, 0, 1),
    @(14, 'Arkansas', 'Missouri', 0, 34, 0, 1),
    @(14, 'Akron', 'Ohio', 0, 11, 0, 1),
    @(14, 'Oklahoma', 'Texas Christian', 0, 24, 1, 0),
    @(14, 'Oregon', 'Oregon State', 0, 24, 1, 0),
    @(14, 'Penn State', 'Michigan State', 1, 42, 1, 0),
    @(14, 'Texas', 'Texas Tech', 0, 50, 1, 0),
    @(14, 'Central Michigan', 'Toledo', 0, 15, 0, 1),
    @(14, 'Tulane', 'Texas-San Antonio', 0, 13, 1, 0),
    @(14, 'New Mexico', 'Utah State', 0, 3, 0, 1),
    @(14, 'Auburn', 'Alabama', 0, 3, 0, 1),
    @(14, 'Appalachian State', 'Georgia Southern', 0, 28, 1, 0),
    @(14, 'Arizona State', 'Arizona', 0, 36, 0, 1),
    @(14, 'UCLA', 'California', 0, 26, 0, 1),
    @(14, 'Central Florida', 'Houston', 0, 14, 1, 0),
    @(14, 'South Carolina', 'Clemson', 0, 9, 0, 1),
    @(14, 'Massachusetts', 'Connecticut', 0, 13, 0, 1),
    @(14, 'Duke', 'Pittsburgh', 0, 11, 1, 0),
    @(14, 'Florida', 'Florida State', 0, 9, 0, 1),
    @(14, 'Georgia Tech', 'Georgia', 0, 8, 0, 1),
    @(14, 'Hawaii', 'Colorado State', 0, 3, 1, 0),
    @(14, 'Kansas State', 'Iowa State', 0, 7, 0, 1),
    @(14, 'Coastal Carolina', 'James Madison', 0, 42, 0, 1),
    @(14, 'Cincinnati', 'Kansas', 0, 33, 0, 1),
    @(14, 'Louisville', 'Kentucky', 0, 7, 0, 1),
    @(14, 'Texas-El Paso', 'Liberty', 0, 14, 0, 1),
    @(14, 'Louisiana', 'Louisiana-Monroe', 0, 31, 1, 0),
    @(14, 'Louisiana State', 'Texas A&M', 0, 12, 1, 0),
    @(14, 'Marshall', 'Arkansas State', 0, 14, 1, 0),
    @(14, 'Rutgers', 'Maryland', 0, 18, 0, 1),
    @(14, 'Ball State', 'Miami (OH)', 0, 2, 0, 1),
    @(14, 'Michigan', 'Ohio State', 0, 6, 1, 0),
    @(14, 'New Mexico State', 'Jacksonville State', 0, 3, 1, 0),
    @(14, 'North Carolina State', 'North Carolina', 0, 19, 1, 0),
    @(14, 'North Texas', 'Alabama-Birmingham', 0, 3, 1, 0),
    @(14, 'Kent State', 'Northern Illinois', 0, 10, 0, 1),
    @(14, 'Illinois', 'Northwestern', 0, 2, 0, 1),
    @(14, 'Stanford', 'Notre Dame', 0, 33, 0, 1),
    @(14, 'Oklahoma State', 'Brigham Young', 0, 6, 1, 0),
    @(14, 'Old Dominion', 'Georgia State', 0, 1, 1, 0),
    @(14, 'Purdue', 'Indiana', 0, 4, 1, 0),
    @(14, 'Rice', 'Florida Atlantic', 0, 3, 1, 0),
    @(14, 'Sam Houston', 'Middle Tennessee State', 0, 3, 1, 0),
    @(14, 'San Diego State', 'Fresno State', 0, 15, 1, 0),
    @(14, 'Nevada-Las Vegas', 'San Jose State', 0, 6, 0, 1),
    @(14, 'South Florida', 'Charlotte', 0, 34, 1, 0),
    @(14, 'Southern Methodist', 'Navy', 0, 45, 1, 0),
    @(14, 'Syracuse', 'Wake Forest', 0, 4, 1, 0),
    @(14, 'Tennessee', 'Vanderbilt', 0, 24, 1, 0),
    @(14, 'Texas State', 'South Alabama', 0, 8, 1, 0),
    @(14, 'Southern Mississippi', 'Troy', 0, 18, 0, 1),
    @(14, 'East Carolina', 'Tulsa', 0, 2, 0, 1),
    @(14, 'Utah', 'Colorado', 0, 6, 1, 0),
    @(14, 'Virginia', 'Virginia Tech', 0, 38, 0, 1),
    @(14, 'Washington', 'Washington State', 0, 3, 1, 0),
    @(14, 'Baylor', 'West Virginia', 0, 3, 0, 1),
    @(14, 'Florida International', 'Western Kentucky', 0, 13, 0, 1),
    @(14, 'Minnesota', 'Wisconsin', 0, 14, 0, 1),
    @(14, 'Nevada', 'Wyoming', 0, 36, 0, 1),
)

$startRow = 794
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
}

Write-Output "Added $($rows.Count) rows starting at $startRow"
